$wb = $excel.ActiveWorkbook

# ---- Sheet1 (BY-Demands): selection change only ----
$wsBY = $wb.Worksheets.Item("BY-Demands")

# ---- Sheet2 (REG_TRA_DEMANDS): selection change, no longer the active tab ----
$wsTra = $wb.Worksheets.Item("REG_TRA_DEMANDS")

# ---- Sheet4 (DEMANDS): content + view changes; becomes the active tab ----
$ws = $wb.Worksheets.Item("DEMANDS")

# New column width for column H (8)
$ws.Columns.Item(8).ColumnWidth = 10.2

# Row 49: new header row (years 2018-2100), matches row 2 header pattern
$ws.Range("F49").Value = "Cset_CN"
$ws.Range("G49").Value = "*Description"
$ws.Range("H49").Value = 2018
$ws.Range("I49").Value = 2019
$ws.Range("J49").Value = 2020
$ws.Range("K49").Value = 2021
$ws.Range("L49").Value = 2022
$ws.Range("M49").Value = 2023
$ws.Range("N49").Value = 2024
$ws.Range("O49").Value = 2025
$ws.Range("P49").Value = 2026
$ws.Range("Q49").Value = 2027
$ws.Range("R49").Value = 2028
$ws.Range("S49").Value = 2029
$ws.Range("T49").Value = 2030
$ws.Range("U49").Value = 2031
$ws.Range("V49").Value = 2032
$ws.Range("W49").Value = 2033
$ws.Range("X49").Value = 2034
$ws.Range("Y49").Value = 2035
$ws.Range("Z49").Value = 2036
$ws.Range("AA49").Value = 2037
$ws.Range("AB49").Value = 2038
$ws.Range("AC49").Value = 2039
$ws.Range("AD49").Value = 2040
$ws.Range("AE49").Value = 2041
$ws.Range("AF49").Value = 2042
$ws.Range("AG49").Value = 2043
$ws.Range("AH49").Value = 2044
$ws.Range("AI49").Value = 2045
$ws.Range("AJ49").Value = 2046
$ws.Range("AK49").Value = 2047
$ws.Range("AL49").Value = 2048
$ws.Range("AM49").Value = 2049
$ws.Range("AN49").Value = 2050
$ws.Range("AO49").Value = 2051
$ws.Range("AP49").Value = 2052
$ws.Range("AQ49").Value = 2053
$ws.Range("AR49").Value = 2054
$ws.Range("AS49").Value = 2055
$ws.Range("AT49").Value = 2056
$ws.Range("AU49").Value = 2057
$ws.Range("AV49").Value = 2058
$ws.Range("AW49").Value = 2059
$ws.Range("AX49").Value = 2060
$ws.Range("AY49").Value = 2061
$ws.Range("AZ49").Value = 2062
$ws.Range("BA49").Value = 2063
$ws.Range("BB49").Value = 2064
$ws.Range("BC49").Value = 2065
$ws.Range("BD49").Value = 2066
$ws.Range("BE49").Value = 2067
$ws.Range("BF49").Value = 2068
$ws.Range("BG49").Value = 2069
$ws.Range("BH49").Value = 2070
$ws.Range("BI49").Value = 2071
$ws.Range("BJ49").Value = 2072
$ws.Range("BK49").Value = 2073
$ws.Range("BL49").Value = 2074
$ws.Range("BM49").Value = 2075
$ws.Range("BN49").Value = 2076
$ws.Range("BO49").Value = 2077
$ws.Range("BP49").Value = 2078
$ws.Range("BQ49").Value = 2079
$ws.Range("BR49").Value = 2080
$ws.Range("BS49").Value = 2081
$ws.Range("BT49").Value = 2082
$ws.Range("BU49").Value = 2083
$ws.Range("BV49").Value = 2084
$ws.Range("BW49").Value = 2085
$ws.Range("BX49").Value = 2086
$ws.Range("BY49").Value = 2087
$ws.Range("BZ49").Value = 2088
$ws.Range("CA49").Value = 2089
$ws.Range("CB49").Value = 2090
$ws.Range("CC49").Value = 2091
$ws.Range("CD49").Value = 2092
$ws.Range("CE49").Value = 2093
$ws.Range("CF49").Value = 2094
$ws.Range("CG49").Value = 2095
$ws.Range("CH49").Value = 2096
$ws.Range("CI49").Value = 2097
$ws.Range("CJ49").Value = 2098
$ws.Range("CK49").Value = 2099
$ws.Range("CL49").Value = 2100

# Rows 50-52: move the original RSD demand values here (verbatim, from old rows 26-28)
$ws.Range("F50").Value = "RSD_Apt"
$ws.Range("G50").Value = "Residential Apartment Demand"
$ws.Range("H50").Value = 213414.55695615499
$ws.Range("I50").Value = 218956.66937259401
$ws.Range("J50").Value = 225300.33816185399
$ws.Range("K50").Value = 242574.41637552099
$ws.Range("L50").Value = 259561.94669652099
$ws.Range("M50").Value = 275023.94138045801
$ws.Range("N50").Value = 289113.78640676098
$ws.Range("O50").Value = 302217.18367944099
$ws.Range("P50").Value = 314549.61337722198
$ws.Range("Q50").Value = 326290.81316684798
$ws.Range("R50").Value = 337584.19698321301
$ws.Range("S50").Value = 348488.670477651
$ws.Range("T50").Value = 359053.362495562
$ws.Range("U50").Value = 379435.10129855498
$ws.Range("V50").Value = 394286.43020683102
$ws.Range("W50").Value = 409460.08200179902
$ws.Range("X50").Value = 424963.52313310798
$ws.Range("Y50").Value = 440776.66278306401
$ws.Range("Z50").Value = 456888.48985664902
$ws.Range("AA50").Value = 473286.23919771903
$ws.Range("AB50").Value = 489960.78944409703
$ws.Range("AC50").Value = 506886.30914768699
$ws.Range("AD50").Value = 524045.59801496298
$ws.Range("AI50").Value = 632981.68201128603
$ws.Range("AN50").Value = 740749.97310115199
$ws.Range("AS50").Value = 853983.36761178297
$ws.Range("AX50").Value = 973313.84240591305
$ws.Range("BC50").Value = 1101920.21948223
$ws.Range("BH50").Value = 1243456.94226919
$ws.Range("BM50").Value = 1400399.1537987201
$ws.Range("BR50").Value = 1574723.2584762699
$ws.Range("BW50").Value = 1766492.4181818999
$ws.Range("CB50").Value = 1974228.0679925601
$ws.Range("CG50").Value = 2195390.7432208802
$ws.Range("CL50").Value = 2428365.4614147102

$ws.Range("F51").Value = "RSD_Att"
$ws.Range("G51").Value = "Residential Attached Demand"
$ws.Range("H51").Value = 766617.95373000402
$ws.Range("I51").Value = 776019.61867223098
$ws.Range("J51").Value = 788300.85620074603
$ws.Range("K51").Value = 801957.14588082698
$ws.Range("L51").Value = 819060.32622000796
$ws.Range("M51").Value = 834661.58347222698
$ws.Range("N51").Value = 848818.67768832797
$ws.Range("O51").Value = 862331.35320143297
$ws.Range("P51").Value = 875509.45116964797
$ws.Range("Q51").Value = 888595.96249835705
$ws.Range("R51").Value = 901764.69015385106
$ws.Range("S51").Value = 914993.375655549
$ws.Range("T51").Value = 928266.37631273095
$ws.Range("U51").Value = 943708.89277841605
$ws.Range("V51").Value = 957570.71779714397
$ws.Range("W51").Value = 971442.33945732203
$ws.Range("X51").Value = 985334.55511306797
$ws.Range("Y51").Value = 999194.66131792602
$ws.Range("Z51").Value = 1012994.63657764
$ws.Range("AA51").Value = 1026705.03546086
$ws.Range("AB51").Value = 1040306.63185868
$ws.Range("AC51").Value = 1053747.0705237
$ws.Range("AD51").Value = 1066995.89427272
$ws.Range("AI51").Value = 1166763.43459146
$ws.Range("AN51").Value = 1244336.6649656901
$ws.Range("AS51").Value = 1314710.96966167
$ws.Range("AX51").Value = 1379878.79642933
$ws.Range("BC51").Value = 1444629.9926038501
$ws.Range("BH51").Value = 1512989.3882673399
$ws.Range("BM51").Value = 1586512.86397531
$ws.Range("BR51").Value = 1665745.9920316101
$ws.Range("BW51").Value = 1749111.94897432
$ws.Range("CB51").Value = 1833896.3857972301
$ws.Range("CG51").Value = 1917030.21036523
$ws.Range("CL51").Value = 1996868.66239624

$ws.Range("F52").Value = "RSD_Det"
$ws.Range("G52").Value = "Residential Detached Demand"
$ws.Range("H52").Value = 721105.90866437601
$ws.Range("I52").Value = 729040.47594293603
$ws.Range("J52").Value = 739683.68570912001
$ws.Range("K52").Value = 748342.44886999705
$ws.Range("L52").Value = 760630.721289968
$ws.Range("M52").Value = 771850.38797134894
$ws.Range("N52").Value = 782013.10939486197
$ws.Range("O52").Value = 791820.47563049395
$ws.Range("P52").Value = 801523.77429367404
$ws.Range("Q52").Value = 811316.73568663199
$ws.Range("R52").Value = 821333.312278467
$ws.Range("S52").Value = 831532.33535063395
$ws.Range("T52").Value = 841882.469887482
$ws.Range("U52").Value = 851415.47129100398
$ws.Range("V52").Value = 861039.36197674705
$ws.Range("W52").Value = 870573.18258039001
$ws.Range("X52").Value = 880025.70362676098
$ws.Range("Y52").Value = 889349.03440936899
$ws.Range("Z52").Value = 898517.81453963905
$ws.Range("AA52").Value = 907505.76600673597
$ws.Range("AB52").Value = 916295.97516395
$ws.Range("AC52").Value = 924842.69231704494
$ws.Range("AD52").Value = 933119.97808990802
$ws.Range("AI52").Value = 1001750.03116784
$ws.Range("AN52").Value = 1047946.93470063
$ws.Range("AS52").Value = 1085051.9905217399
$ws.Range("AX52").Value = 1114915.32706986
$ws.Range("BC52").Value = 1141469.10325923
$ws.Range("BH52").Value = 1167711.9675762199
$ws.Range("BM52").Value = 1194473.05229634
$ws.Range("BR52").Value = 1221698.7321783099
$ws.Range("BW52").Value = 1247749.8188300601
$ws.Range("CB52").Value = 1270297.7134954601
$ws.Range("CG52").Value = 1286979.0355445601
$ws.Range("CL52").Value = 1296606.1476706499

# Rows 26-28: replace literal values with formulas referencing new rows 50-52 (convert to original units, /1000)
$ws.Range("H26:AD26").Formula = "=H50/1000"
$ws.Range("AI26").Formula = "=AI50/1000"
$ws.Range("AN26").Formula = "=AN50/1000"
$ws.Range("AS26").Formula = "=AS50/1000"
$ws.Range("AX26").Formula = "=AX50/1000"
$ws.Range("BC26").Formula = "=BC50/1000"
$ws.Range("BH26").Formula = "=BH50/1000"
$ws.Range("BM26").Formula = "=BM50/1000"
$ws.Range("BR26").Formula = "=BR50/1000"
$ws.Range("BW26").Formula = "=BW50/1000"
$ws.Range("CB26").Formula = "=CB50/1000"
$ws.Range("CG26").Formula = "=CG50/1000"
$ws.Range("CL26").Formula = "=CL50/1000"

$ws.Range("H27:AD27").Formula = "=H51/1000"
$ws.Range("AI27").Formula = "=AI51/1000"
$ws.Range("AN27").Formula = "=AN51/1000"
$ws.Range("AS27").Formula = "=AS51/1000"
$ws.Range("AX27").Formula = "=AX51/1000"
$ws.Range("BC27").Formula = "=BC51/1000"
$ws.Range("BH27").Formula = "=BH51/1000"
$ws.Range("BM27").Formula = "=BM51/1000"
$ws.Range("BR27").Formula = "=BR51/1000"
$ws.Range("BW27").Formula = "=BW51/1000"
$ws.Range("CB27").Formula = "=CB51/1000"
$ws.Range("CG27").Formula = "=CG51/1000"
$ws.Range("CL27").Formula = "=CL51/1000"

$ws.Range("H28:AD28").Formula = "=H52/1000"
$ws.Range("AI28").Formula = "=AI52/1000"
$ws.Range("AN28").Formula = "=AN52/1000"
$ws.Range("AS28").Formula = "=AS52/1000"
$ws.Range("AX28").Formula = "=AX52/1000"
$ws.Range("BC28").Formula = "=BC52/1000"
$ws.Range("BH28").Formula = "=BH52/1000"
$ws.Range("BM28").Formula = "=BM52/1000"
$ws.Range("BR28").Formula = "=BR52/1000"
$ws.Range("BW28").Formula = "=BW52/1000"
$ws.Range("CB28").Formula = "=CB52/1000"
$ws.Range("CG28").Formula = "=CG52/1000"
$ws.Range("CL28").Formula = "=CL52/1000"

# Row 40: new I40 zero cell
$ws.Range("I40").Value = 0

# ---- Selections / active sheet (applied last so the final selection state sticks) ----
$wsBY.Activate()
$wsBY.Range("A14").Select()

$wsTra.Activate()
$wsTra.Range("G84").Select()

$ws.Activate()
$ws.Range("CL28").Select()
